$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New reference/predictor measurement rows added to each of the three
# "Instance N" result tables (rows 43-46, 57-60, 71-74). Columns B/C/D are
# the raw Run 1/2/3 T_inf readings, L/M/N are the raw ambient-temperature
# readings; columns E-K are pre-existing shared formulas that recompute
# automatically once the raw inputs below are populated.
$data = @(
    @{Row=43; A="Mod.2 B predictor(2)";     B=54.591999999999999; C=53.94;               D=52.802;             L=23.838999999999999; M=24.271999999999998; N=23.225200000000001},
    @{Row=44; A="Mod.2 B predictor(2) LTF"; B=54.648000000000003; C=54.067999999999998;  D=52.91;              L=23.662599999999902; M=24.376999999999999; N=23.165199999999999},
    @{Row=45; A="ilp reference";            B=54.723999999999997; C=54.835999999999999;  D=53.276000000000003;L=23.587599999999998; M=25.370399999999901; N=23.058},
    @{Row=46; A="reference";                B=55.643999999999998; C=55.911999999999999;  D=54.04;              L=23.523599999999998; M=25.434799999999999; N=23.034399999999899},

    @{Row=57; A="Mod.2 B predictor(2)";     B=54.067999999999998; C=54.612000000000002;  D=52.561999999999998;L=23.471800000000002; M=25.880800000000001; N=22.965599999999998},
    @{Row=58; A="Mod.2 B predictor(2) LTF"; B=55.735999999999997; C=57.584000000000003;  D=54.756;             L=23.413999999999898; M=26.7837999999999;   N=22.918199999999999},
    @{Row=59; A="ilp reference";            B=54.847999999999999; C=55.508000000000003;  D=53.332000000000001;L=23.324000000000002; M=25.121399999999898; N=22.819599999999902},
    @{Row=60; A="reference";                B=56.96;              C=57.868000000000002;  D=55.968000000000004;L=23.2898;             M=25.0654;            N=22.8216},

    @{Row=71; A="Mod.2 B predictor(2)";     B=52.722000000000001; C=54.508000000000003;  D=52.415999999999997;L=23.240199999999898; M=24.883199999999999; N=22.734000000000002},
    @{Row=72; A="Mod.2 B predictor(2) LTF"; B=53.7;               C=55.648000000000003;  D=53.031999999999996;L=23.175799999999999; M=25.288599999999999; N=22.718599999999999},
    @{Row=73; A="ilp reference";            B=54.223999999999997; C=55.643999999999998;  D=53.692;             L=23.138999999999999; M=24.661999999999999; N=22.7362},
    @{Row=74; A="reference";                B=53.704000000000001; C=55.308;               D=53.433999999999997;L=23.128799999999998; M=24.349;             N=22.787399999999899}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A    # A: label
    $ws.Cells.Item($r, 2).Value = $item.B     # B: Run 1
    $ws.Cells.Item($r, 3).Value = $item.C     # C: Run 2
    $ws.Cells.Item($r, 4).Value = $item.D     # D: Run 3
    $ws.Cells.Item($r, 12).Value = $item.L    # L: Ambient Run 1
    $ws.Cells.Item($r, 13).Value = $item.M    # M: Ambient Run 2
    $ws.Cells.Item($r, 14).Value = $item.N    # N: Ambient Run 3
}

# Selection moved to the newly-populated ambient-temperature cells of the
# "Instance 6" table.
$ws.Range("N71:N74").Select()
